$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings are not
# auto-converted to numbers by the COM Value setter (matches source data,
# which stores these as literal text).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.81'
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.44%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '3.37'
$ws.Range("E38").Value = '  -4.88%  '

$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.435.18'
$ws.Range("E41").Value = '  +1.98%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.16'
$ws.Range("E42").Value = '  -0.80%  '

$ws.Range("D2").Value = '37.208.33'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '1.994.30'
$ws.Range("E3").Value = '  -1.74%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").Value = '254.36'
$ws.Range("E5").Value = '  +3.32%  '
$ws.Range("D6").Value = '0.606'
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").Value = '54.90'
$ws.Range("E8").Value = '  -6.06%  '
$ws.Range("D9").Value = '0.376'
$ws.Range("E9").Value = '  -3.27%  '
$ws.Range("D10").Value = '0.0758'
$ws.Range("E10").Value = '  -5.54%  '
$ws.Range("D11").Value = '0.0998'
$ws.Range("E11").Value = '  -3.46%  '
$ws.Range("D12").Value = '2.301.09'
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").Value = '14.02'
$ws.Range("E13").Value = '  -6.10%  '
$ws.Range("D14").Value = '21.40'
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").Value = '0.765'
$ws.Range("E15").Value = '  -8.10%  '
$ws.Range("D16").Value = '5.13'
$ws.Range("E16").Value = '  -4.70%  '
$ws.Range("D17").Value = '1.990.67'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("D18").Value = '37.129.49'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '69.45'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").Value = '0.0₃0821'
$ws.Range("E20").Value = '  -4.14%  '
$ws.Range("D21").Value = '232.60'
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").Value = '5.02'
$ws.Range("E22").Value = '  -4.22%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = '2.50'
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").Value = '164.34'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").Value = '8.74'
$ws.Range("E27").Value = '  -5.41%  '
$ws.Range("E28").Value = '  -5.03%  '
$ws.Range("D29").Value = '19.27'
$ws.Range("E29").Value = '  -3.51%  '
$ws.Range("D30").Value = '1.29'
$ws.Range("E30").Value = '  -5.17%  '
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").Value = '4.51'
$ws.Range("E32").Value = '  -5.34%  '
$ws.Range("D33").Value = '0.0619'
$ws.Range("E33").Value = '  -7.71%  '
$ws.Range("D34").Value = '4.32'
$ws.Range("E34").Value = '  -4.43%  '
$ws.Range("D35").Value = '2.36'
$ws.Range("E35").Value = '  -4.67%  '
$ws.Range("D39").Value = '5.30'
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("D40").Value = '3.06'
$ws.Range("E40").Value = '  +2.08%  '
$ws.Range("D43").Value = '0.0906'
$ws.Range("E43").Value = '  -6.94%  '
$ws.Range("D44").Value = '0.0206'
$ws.Range("E44").Value = '  -4.80%  '
$ws.Range("D45").Value = '15.62'
$ws.Range("E45").Value = '  -5.45%  '
$ws.Range("D46").Value = '88.41'
$ws.Range("E46").Value = '  -3.18%  '
$ws.Range("D47").Value = '1.01'
$ws.Range("E47").Value = '  -3.23%  '
$ws.Range("D48").Value = '2.89'
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("D49").Value = '6.73'
$ws.Range("E49").Value = '  -9.97%  '
$ws.Range("D50").Value = '2.190.56'
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '1.88'
$ws.Range("E51").Value = '  -9.64%  '

# Restore default styling (remove the temporary Text number format) so
# cells keep their original (unstyled) appearance while staying text.
$ws.Range("D2:E51").Style = "Normal"
